# This script applies a re-shuffling of the weekly price records (rows 2-30)
# in the "Fruta, Agrícola del Norte S.A. de Arica - Ciruela" sheet, matching
# the new weekly snapshot of market data. Columns A, B, C, E, F, G, H, I, J
# (market/product identity columns) stay put; columns D, K, L, M, N, O, P, Q,
# R, S, T (date, variety, quality, volume, prices, unit, origin, $/Kg, Kg/unit)
# are rearranged among the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) data for columns D..T on every data row,
# before any writes happen, so the permutation below is based on consistent
# "before" values regardless of write order.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $snapshot[$r] = $ws.Range("D$r`:T$r").Value2
}

# Re-distribute the snapshot rows to their new positions.

$ws.Range("D2:T2").Value = $snapshot[6]
$ws.Range("D3:T3").Value = $snapshot[8]
$ws.Range("D4:T4").Value = $snapshot[20]
$ws.Range("D5:T5").Value = $snapshot[29]
$ws.Range("D6:T6").Value = $snapshot[16]
$ws.Range("D7:T7").Value = $snapshot[30]
$ws.Range("D8:T8").Value = $snapshot[7]
$ws.Range("D9:T9").Value = $snapshot[21]
$ws.Range("D10:T10").Value = $snapshot[17]
$ws.Range("D11:T11").Value = $snapshot[18]
$ws.Range("D12:T12").Value = $snapshot[14]
$ws.Range("D13:T13").Value = $snapshot[2]
$ws.Range("D14:T14").Value = $snapshot[26]
$ws.Range("D15:T15").Value = $snapshot[27]
$ws.Range("D16:T16").Value = $snapshot[4]
$ws.Range("D17:T17").Value = $snapshot[11]
$ws.Range("D18:T18").Value = $snapshot[13]
$ws.Range("D19:T19").Value = $snapshot[12]
$ws.Range("D20:T20").Value = $snapshot[24]
$ws.Range("D21:T21").Value = $snapshot[23]
$ws.Range("D22:T22").Value = $snapshot[10]
$ws.Range("D23:T23").Value = $snapshot[25]
$ws.Range("D24:T24").Value = $snapshot[22]
$ws.Range("D25:T25").Value = $snapshot[15]
$ws.Range("D26:T26").Value = $snapshot[9]
$ws.Range("D27:T27").Value = $snapshot[5]
$ws.Range("D28:T28").Value = $snapshot[19]
$ws.Range("D29:T29").Value = $snapshot[28]
$ws.Range("D30:T30").Value = $snapshot[3]
